# Auto update stock data
# Updates the Date_1 column to 2026/01/13 and refreshes EBITDA (and for one
# row, the Debt / Equity Ratio) figures for the latest observation of each
# company block in the risk-scores sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    # Force the cell to stay a text value (matches the workbook's existing
    # inline-string convention for this column) instead of letting Excel
    # auto-coerce date-looking / numeric-looking strings.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
}

$updates = @(
    @{ Row = 2;  A = "2026/01/13"; B = "8.21" }
    @{ Row = 8;  A = "2026/01/13"; B = "8.62" }
    @{ Row = 14; A = "2026/01/13"; B = "3.15" }
    @{ Row = 20; A = "2026/01/13"; B = "13.54" }
    @{ Row = 26; A = "2026/01/13"; B = "11.87" }
    @{ Row = 32; A = "2026/01/13"; B = "28.95" }
    @{ Row = 38; A = "2026/01/13" }
    @{ Row = 44; A = "2026/01/13"; B = "14.32" }
    @{ Row = 50; A = "2026/01/13"; B = "11.87" }
    @{ Row = 56; A = "2026/01/13"; B = "32.96" }
    @{ Row = 62; A = "2026/01/13"; B = "11.46"; C = "2.56" }
    @{ Row = 68; A = "2026/01/13"; B = "12.95" }
    @{ Row = 74; A = "2026/01/13"; B = "18.40" }
)

foreach ($u in $updates) {
    Set-TextValue $ws.Cells.Item($u.Row, 1) $u.A
    if ($u.ContainsKey("B")) {
        Set-TextValue $ws.Cells.Item($u.Row, 2) $u.B
    }
    if ($u.ContainsKey("C")) {
        Set-TextValue $ws.Cells.Item($u.Row, 3) $u.C
    }
}
